# Final re-run of "Impact of the psychosis disorder on education"
# (N=146 -> N=143), including fig: unmet social support.
#
# All plotted shapes (bars, gridlines, tick labels, data labels,
# title) live inside one top-level group shape on slide 1; we
# reach them by name via GroupItems so the script does not rely
# on brittle numeric indices.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(1)

# pl5 (id=5)
$sh = $grp.GroupItems.Item('pl5')
$sh.Top = 363.73126220703125

# pl6 (id=6)
$sh = $grp.GroupItems.Item('pl6')
$sh.Top = 248.02740478515625

# pl7 (id=7)
$sh = $grp.GroupItems.Item('pl7')
$sh.Top = 132.32354736328125

# rc8 (id=8)
$sh = $grp.GroupItems.Item('rc8')
$sh.Top = 455.16168212890625
$sh.Height = 24.273466110229492

# rc9 (id=9)
$sh = $grp.GroupItems.Item('rc9')
$sh.Top = 228.6085968017578
$sh.Height = 250.8265380859375

# rc10 (id=10)
$sh = $grp.GroupItems.Item('rc10')
$sh.Top = 244.79095458984375
$sh.Height = 234.64418029785156

# rc11 (id=11)
$sh = $grp.GroupItems.Item('rc11')
$sh.Top = 422.79693603515625
$sh.Height = 56.63819122314453

# rc12 (id=12)
$sh = $grp.GroupItems.Item('rc12')
$sh.Top = 325.7027587890625
$sh.Height = 153.7323760986328

# rc13 (id=13)
$sh = $grp.GroupItems.Item('rc13')
$sh.Top = 333.7939453125
$sh.Height = 145.6411895751953

# tx15 (id=15)
$sh = $grp.GroupItems.Item('tx15')
$sh.Top = 408.93505859375
$sh.Height = 10.40574836730957
$sh.TextFrame.TextRange.Text = '3'

# tx16 (id=16)
$sh = $grp.GroupItems.Item('tx16')
$sh.Top = 426.47576904296875
$sh.TextFrame.TextRange.Text = '(2%)'

# tx17 (id=17)
$sh = $grp.GroupItems.Item('tx17')
$sh.Top = 182.38204956054688

# tx18 (id=18)
$sh = $grp.GroupItems.Item('tx18')
$sh.Top = 199.92276000976562
$sh.TextFrame.TextRange.Text = '(22%)'

# tx19 (id=19)
$sh = $grp.GroupItems.Item('tx19')
$sh.Top = 198.57135009765625
$sh.Height = 10.398818969726562
$sh.TextFrame.TextRange.Text = '29'

# tx20 (id=20)
$sh = $grp.GroupItems.Item('tx20')
$sh.Top = 216.10513305664062
$sh.TextFrame.TextRange.Text = '(20%)'

# tx21 (id=21)
$sh = $grp.GroupItems.Item('tx21')
$sh.Top = 376.9246520996094

# tx22 (id=22)
$sh = $grp.GroupItems.Item('tx22')
$sh.Top = 394.1110534667969

# tx23 (id=23)
$sh = $grp.GroupItems.Item('tx23')
$sh.Top = 279.483154296875

# tx24 (id=24)
$sh = $grp.GroupItems.Item('tx24')
$sh.Top = 297.0168762207031

# tx25 (id=25)
$sh = $grp.GroupItems.Item('tx25')
$sh.Top = 287.5743408203125
$sh.TextFrame.TextRange.Text = '18'

# tx26 (id=26)
$sh = $grp.GroupItems.Item('tx26')
$sh.Top = 305.1080322265625

# tx27 (id=27)
$sh = $grp.GroupItems.Item('tx27')
$sh.TextFrame.TextRange.Text = '36'

# tx31 (id=31)
$sh = $grp.GroupItems.Item('tx31')
$sh.Top = 359.2547302246094

# tx32 (id=32)
$sh = $grp.GroupItems.Item('tx32')
$sh.Top = 243.55087280273438

# tx33 (id=33)
$sh = $grp.GroupItems.Item('tx33')
$sh.Top = 127.8411865234375

# pl35 (id=35)
$sh = $grp.GroupItems.Item('pl35')
$sh.Top = 363.73126220703125

# pl36 (id=36)
$sh = $grp.GroupItems.Item('pl36')
$sh.Top = 248.02740478515625

# pl37 (id=37)
$sh = $grp.GroupItems.Item('pl37')
$sh.Top = 132.32354736328125

# tx54 (id=54)
$sh = $grp.GroupItems.Item('tx54')
$sh.TextFrame.TextRange.Text = 'Impact of the psychosis disorder on education (N=143)'
